$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (MBOM instead of Sheet1); defined names follow automatically
$ws.Name = "MBOM"

# Remove the now-stale trailing rows (92-99) - new data only runs through row 91
$ws.Rows("92:99").Delete()

# Update Name (B) / Q (C) / Referencia (D) columns for every data row with the new
# expression-based values (Refuerzo1..Refuerzo90 + re-derived Referencia lookups)
$ws.Range("B2").Value = "Refuerzo1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "_110"
$ws.Range("B3").Value = "Refuerzo2"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "_117"
$ws.Range("B4").Value = "Refuerzo3"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "_125"
$ws.Range("B5").Value = "Refuerzo4"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "_111"
$ws.Range("B6").Value = "Refuerzo5"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "_119"
$ws.Range("B7").Value = "Refuerzo6"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "_107"
$ws.Range("B8").Value = "Refuerzo7"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "_108"
$ws.Range("B9").Value = "Refuerzo8"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "_108"
$ws.Range("B10").Value = "Refuerzo9"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "_103"
$ws.Range("B11").Value = "Refuerzo10"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "_118"
$ws.Range("B12").Value = "Refuerzo11"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "_101"
$ws.Range("B13").Value = "Refuerzo12"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "_123"
$ws.Range("B14").Value = "Refuerzo13"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "_125"
$ws.Range("B15").Value = "Refuerzo14"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "_122"
$ws.Range("B16").Value = "Refuerzo15"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "_116"
$ws.Range("B17").Value = "Refuerzo16"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "_105"
$ws.Range("B18").Value = "Refuerzo17"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "_104"
$ws.Range("B19").Value = "Refuerzo18"
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = "_119"
$ws.Range("B20").Value = "Refuerzo19"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "_106"
$ws.Range("B21").Value = "Refuerzo20"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "_112"
$ws.Range("B22").Value = "Refuerzo21"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "_104"
$ws.Range("B23").Value = "Refuerzo22"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "_116"
$ws.Range("B24").Value = "Refuerzo23"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = "_120"
$ws.Range("B25").Value = "Refuerzo24"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "_104"
$ws.Range("B26").Value = "Refuerzo25"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "_101"
$ws.Range("B27").Value = "Refuerzo26"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "_100"
$ws.Range("B28").Value = "Refuerzo27"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "_105"
$ws.Range("B29").Value = "Refuerzo28"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "_108"
$ws.Range("B30").Value = "Refuerzo29"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "_122"
$ws.Range("B31").Value = "Refuerzo30"
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = "_126"
$ws.Range("B32").Value = "Refuerzo31"
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = "_124"
$ws.Range("B33").Value = "Refuerzo32"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = "_129"
$ws.Range("B34").Value = "Refuerzo33"
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = "_114"
$ws.Range("B35").Value = "Refuerzo34"
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = "_128"
$ws.Range("B36").Value = "Refuerzo35"
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = "_120"
$ws.Range("B37").Value = "Refuerzo36"
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = "_103"
$ws.Range("B38").Value = "Refuerzo37"
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = "_137"
$ws.Range("B39").Value = "Refuerzo38"
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = "_110"
$ws.Range("B40").Value = "Refuerzo39"
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = "_103"
$ws.Range("B41").Value = "Refuerzo40"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = "_111"
$ws.Range("B42").Value = "Refuerzo41"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = "_125"
$ws.Range("B43").Value = "Refuerzo42"
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = "_123"
$ws.Range("B44").Value = "Refuerzo43"
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = "_129"
$ws.Range("B45").Value = "Refuerzo44"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = "_102"
$ws.Range("B46").Value = "Refuerzo45"
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = "_109"
$ws.Range("B47").Value = "Refuerzo46"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = "_113"
$ws.Range("B48").Value = "Refuerzo47"
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = "_114"
$ws.Range("B49").Value = "Refuerzo48"
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = "_122"
$ws.Range("B50").Value = "Refuerzo49"
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = "_105"
$ws.Range("B51").Value = "Refuerzo50"
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = "_113"
$ws.Range("B52").Value = "Refuerzo51"
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = "_116"
$ws.Range("B53").Value = "Refuerzo52"
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = "_100"
$ws.Range("B54").Value = "Refuerzo53"
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = "_117"
$ws.Range("B55").Value = "Refuerzo54"
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = "_137"
$ws.Range("B56").Value = "Refuerzo55"
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = "_111"
$ws.Range("B57").Value = "Refuerzo56"
$ws.Range("C57").Value = 1
$ws.Range("D57").Value = "_124"
$ws.Range("B58").Value = "Refuerzo57"
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = "_115"
$ws.Range("B59").Value = "Refuerzo58"
$ws.Range("C59").Value = 1
$ws.Range("D59").Value = "_126"
$ws.Range("B60").Value = "Refuerzo59"
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = "_101"
$ws.Range("B61").Value = "Refuerzo60"
$ws.Range("C61").Value = 1
$ws.Range("D61").Value = "_102"
$ws.Range("B62").Value = "Refuerzo61"
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = "_121"
$ws.Range("B63").Value = "Refuerzo62"
$ws.Range("C63").Value = 1
$ws.Range("D63").Value = "_117"
$ws.Range("B64").Value = "Refuerzo63"
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = "_107"
$ws.Range("B65").Value = "Refuerzo64"
$ws.Range("C65").Value = 1
$ws.Range("D65").Value = "_106"
$ws.Range("B66").Value = "Refuerzo65"
$ws.Range("C66").Value = 1
$ws.Range("D66").Value = "_126"
$ws.Range("B67").Value = "Refuerzo66"
$ws.Range("C67").Value = 1
$ws.Range("D67").Value = "_118"
$ws.Range("B68").Value = "Refuerzo67"
$ws.Range("C68").Value = 1
$ws.Range("D68").Value = "_123"
$ws.Range("B69").Value = "Refuerzo68"
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = "_113"
$ws.Range("B70").Value = "Refuerzo69"
$ws.Range("C70").Value = 1
$ws.Range("D70").Value = "_120"
$ws.Range("B71").Value = "Refuerzo70"
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = "_115"
$ws.Range("B72").Value = "Refuerzo71"
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = "_112"
$ws.Range("B73").Value = "Refuerzo72"
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = "_115"
$ws.Range("B74").Value = "Refuerzo73"
$ws.Range("C74").Value = 1
$ws.Range("D74").Value = "_110"
$ws.Range("B75").Value = "Refuerzo74"
$ws.Range("C75").Value = 1
$ws.Range("D75").Value = "_121"
$ws.Range("B76").Value = "Refuerzo75"
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = "_107"
$ws.Range("B77").Value = "Refuerzo76"
$ws.Range("C77").Value = 1
$ws.Range("D77").Value = "_109"
$ws.Range("B78").Value = "Refuerzo77"
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = "_129"
$ws.Range("B79").Value = "Refuerzo78"
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = "_128"
$ws.Range("B80").Value = "Refuerzo79"
$ws.Range("C80").Value = 1
$ws.Range("D80").Value = "_118"
$ws.Range("B81").Value = "Refuerzo80"
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = "_100"
$ws.Range("B82").Value = "Refuerzo81"
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = "_128"
$ws.Range("B83").Value = "Refuerzo82"
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = "_112"
$ws.Range("B84").Value = "Refuerzo83"
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = "_102"
$ws.Range("B85").Value = "Refuerzo84"
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = "_137"
$ws.Range("B86").Value = "Refuerzo85"
$ws.Range("C86").Value = 2
$ws.Range("D86").Value = "_119"
$ws.Range("B87").Value = "Refuerzo86"
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = "_109"
$ws.Range("B88").Value = "Refuerzo87"
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = "_114"
$ws.Range("B89").Value = "Refuerzo88"
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = "_106"
$ws.Range("B90").Value = "Refuerzo89"
$ws.Range("C90").Value = 1
$ws.Range("D90").Value = "_124"
$ws.Range("B91").Value = "Refuerzo90"
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = "_121"
# Remove AutoFilter (no longer needed on the cleaned-up sheet)
if ($ws.AutoFilterMode) { $ws.AutoFilterMode = $false }

# Re-fit column B to the new (longer) Name values
$ws.Columns("B:B").AutoFit()
